# Apply updates to the daily/intraday volume table (Table 2) for the
# Bond Futures (row 2/3) and E-mini Futures (row 26/27/28) sections,
# per the commit "Update daily and intraday volume tables for Bond and
# E-mini Futures".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Panel A / FF1, "Ann Window Volume") ---
$ws.Range("K2").Value = 67.41014044196297
$ws.Range("W2").Value = 70.99128985024717

# --- Row 3 (Panel A / FF1, "Diff (Ann - Non)") ---
$ws.Range("J3").Value = 20.17086733588503
$ws.Range("P3").Value = 31.99000618924527

# --- Row 26 (Panel B / Emini, "Ann Window Volume") ---
$ws.Range("D26").Value  = 2550.204933586338
$ws.Range("E26").Value  = 2478.414796679735
$ws.Range("F26").Value  = 433.3548387096774
$ws.Range("G26").Value  = 1978.806451612903
$ws.Range("H26").Value  = 4081.548387096774
$ws.Range("I26").Value  = 221
$ws.Range("J26").Value  = 2862.071656405311
$ws.Range("K26").Value  = 2260.661787173694
$ws.Range("L26").Value  = 1179.360655737705
$ws.Range("M26").Value  = 2757.934426229508
$ws.Range("N26").Value  = 4255.55737704918
$ws.Range("O26").Value  = 221
$ws.Range("P26").Value  = 2870.764219737482
$ws.Range("Q26").Value  = 2183.145465792814
$ws.Range("R26").Value  = 1476.380165289256
$ws.Range("S26").Value  = 2844.611570247934
$ws.Range("T26").Value  = 4258.074380165289
$ws.Range("U26").Value  = 221
$ws.Range("V26").Value  = 2563.625813402284
$ws.Range("W26").Value  = 1910.955060900347
$ws.Range("X26").Value  = 1270.509523809524
$ws.Range("Y26").Value  = 2524.238095238095
$ws.Range("Z26").Value  = 3812.452380952381
$ws.Range("AA26").Value = 221
$ws.Range("AB26").Value = 736.1477546962842
$ws.Range("AC26").Value = 575.0646640001602

# --- Row 27 (Panel B / Emini, "Diff (Ann - Non)") ---
$ws.Range("D27").Value  = 935.2473361553059
$ws.Range("J27").Value  = 1222.043950745493
$ws.Range("P27").Value  = 1220.935361430014
$ws.Range("V27").Value  = 839.0911387631976
$ws.Range("AB27").Value = 95.65479055258466

# --- Row 28 (Panel B / Emini, "# Obs") ---
$ws.Range("D28").Value = 221
$ws.Range("J28").Value = 221
$ws.Range("P28").Value = 221
$ws.Range("V28").Value = 221
